$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table - R_Article")
$ws.Rows.Item(10).Insert()
$ws.Range("A10:F10").Borders.Item(9).LineStyle = 1
$r = $ws.Range("A10")
Write-Host ("A10.Value()=" + $r.Value())
